$wb = $excel.ActiveWorkbook

# --- Fix misspelled shared string ---
# Sheet "AGE" (sheet1), cell E1 held the misspelled shared string
# "Feed convertion ratio, kg". Correct it to "Feed conversion ratio, kg".
# Because the other sheets (CUMULATIVE / P-VALUES / P_cumulative) already
# reference the correctly spelled shared string, fixing this cell makes the
# misspelled shared-string entry unused; on save the engine drops it from
# the shared-strings table and every other sheet's reference re-indexes
# automatically (matches uniqueCount 7 -> 6 in the diff).
$wsAge = $wb.Worksheets.Item("AGE")
$wsAge.Range("E1").Value = "Feed conversion ratio, kg"

# --- Column E width on AGE sheet ---
# Widen column E (bestFit width grew from 19.140625 to 23.28515625 chars).
$wsAge.Columns.Item(5).ColumnWidth = 22.5

# --- Update selections on the other sheets (also updates which tab is
#     marked as selected, since selecting a range activates that sheet) ---

# CUMULATIVE (sheet2): selection -> B1:D1
$wsCum = $wb.Worksheets.Item("CUMULATIVE")
$null = $wsCum.Range("B1:D1").Select()

# P_cumulative (sheet4): selection -> A1:C1
$wsPCum = $wb.Worksheets.Item("P_cumulative")
$null = $wsPCum.Range("A1:C1").Select()

# P-VALUES (sheet3): selection -> B1:D1 (this sheet loses the "tabSelected"
# flag in the target file, so it must not be the last sheet selected)
$wsPVal = $wb.Worksheets.Item("P-VALUES")
$null = $wsPVal.Range("B1:D1").Select()

# AGE (sheet1): selection -> F1. Select this last so AGE ends up the active
# / tabSelected sheet (matches activeTab moving from P-VALUES to AGE).
$null = $wsAge.Range("F1").Select()
